# Daily attendance processing - 2025-10-27 11:43:05
# Normalize the "Recorded By" (column G) entries so that "System" is
# listed first, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# and the two-admin case "dnasr281@gmail.com, admin@admin.com" ->
# "admin@admin.com, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "admin@admin.com, System") {
        $cell.Value2 = "System, admin@admin.com"
    }
    elseif ($val -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, dnasr281@gmail.com"
    }
}
